$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the rolled ability score values into column B for each ability row
$ws.Range("B2").Value = 14
$ws.Range("B3").Value = 6
$ws.Range("B4").Value = 11
$ws.Range("B5").Value = 13
$ws.Range("B6").Value = 15
$ws.Range("B7").Value = 12

# Update the active selection to match where the user ended up
$ws.Range("C16").Select()

$wb.Save()
